$d = $word.ActiveDocument

# 1) Insert "Langages : python, matlab, c, c++" as a new paragraph right
#    before "Bases de données : ..." (it used to be the very last skill
#    line; it is being moved to the top of the list).
$rng = $d.Content
$rng.Find.Execute("Bases de données : SQL, MongoDB, Neo4j, Redis", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint = $d.Range($rng.Start, $rng.Start)
$insertPoint.InsertBefore("Langages : python, matlab, c, c++`r")

# 2) Insert the new "Data Science : ..." paragraph right after
#    "Bases de données : ..." (i.e. right before "Visualisation : tableau").
$rng2 = $d.Content
$rng2.Find.Execute("Visualisation : tableau", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$insertPoint2 = $d.Range($rng2.Start, $rng2.Start)
$insertPoint2.InsertBefore("Data Science : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn`r")

# 3) Rename "MLOps" category to "Machine Learning".
$d.Content.Find.Execute("MLOps : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", $true, $false, $false, $false, $false, $true, 1, $false, "Machine Learning : Git, DVC, Flask, Docker, Github Actions, Heroku, MLflow, Streamlit", 2) | Out-Null

# 4) Replace the old "ML/AI : ..." line (now redundant with the new
#    "Data Science" line above) with the new "Autres" category.
$d.Content.Find.Execute("ML/AI : Scikit-Learn, Keras, Tensorflow, Pandas, pySpark, XGboost, OpenCV, Matplotlib, Seaborn", $true, $false, $false, $false, $false, $true, 1, $false, "Autres : securite, bonnes pratiques de code", 2) | Out-Null

# 5) Remove the old trailing "Langages : ..." paragraph (it was moved to
#    the top of the skills list in step 1). Two occurrences of that text
#    now exist in the document (the new one inserted at the top of the
#    skills block, and the original one at the bottom) - skip past the
#    first one and delete the paragraph containing the second.
$rng3 = $d.Content
$rng3.Find.Execute("Langages : python, matlab, c, c++", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Start = $rng3.End
$rng3.End = $d.Content.End
$rng3.Find.Execute("Langages : python, matlab, c, c++", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oldLangagesPara = $rng3.Paragraphs(1)
$oldLangagesPara.Range.Delete()
